# Rename the existing sheet to "Sales" and add a new "Summary" sheet after it.
$wb = $excel.ActiveWorkbook
$sales = $wb.ActiveSheet
$sales.Name = "Sales"

$summary = $wb.Worksheets.Add($null, $sales)
$summary.Name = "Summary"

# Match the page-margin defaults used on the Sales sheet (0.75/1/0.5 in).
$summary.PageSetup.LeftMargin = 54
$summary.PageSetup.RightMargin = 54
$summary.PageSetup.TopMargin = 72
$summary.PageSetup.BottomMargin = 72
$summary.PageSetup.HeaderMargin = 36
$summary.PageSetup.FooterMargin = 36

# ---- Sales sheet ----
$salesHeaders = @("Region", "Product", "Q1", "Q2", "Q3", "Q4")
for ($i = 0; $i -lt $salesHeaders.Length; $i++) {
    $cell = $sales.Cells.Item(1, $i + 1)
    $cell.Value = $salesHeaders[$i]
    $cell.Font.Bold = $true
}

$salesRows = @(
    @("North", "Widget A", 120, 135, 148, 162),
    @("North", "Widget B", 85, 92, 88, 95),
    @("South", "Widget A", 200, 210, 195, 220),
    @("South", "Widget B", 150, 165, 172, 180),
    @("East", "Gadget C", 75, 80, 82, 90),
    @("West", "Gadget C", 60, 65, 70, 78)
)

for ($r = 0; $r -lt $salesRows.Length; $r++) {
    $row = $salesRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sales.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---- Summary sheet ----
$summaryHeaders = @("Metric", "Value")
for ($i = 0; $i -lt $summaryHeaders.Length; $i++) {
    $cell = $summary.Cells.Item(1, $i + 1)
    $cell.Value = $summaryHeaders[$i]
    $cell.Font.Bold = $true
}

$summary.Cells.Item(2, 1).Value = "Total Q1"
$summary.Cells.Item(2, 2).Formula = "=SUM(Sales!C2:C7)"

$summary.Cells.Item(3, 1).Value = "Total Q2"
$summary.Cells.Item(3, 2).Formula = "=SUM(Sales!D2:D7)"

$summary.Cells.Item(4, 1).Value = "Total Q3"
$summary.Cells.Item(4, 2).Formula = "=SUM(Sales!E2:E7)"

$summary.Cells.Item(5, 1).Value = "Total Q4"
$summary.Cells.Item(5, 2).Formula = "=SUM(Sales!F2:F7)"

$summary.Cells.Item(6, 1).Value = "Grand Total"
$summary.Cells.Item(6, 2).Formula = "=SUM(B2:B5)"
